$d = $word.ActiveDocument

function Get-ParaText($p) {
    return $p.Range.Text.TrimEnd([char]13, [char]7)
}

function Find-ParaIndex($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ((Get-ParaText $p) -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-ParaText($doc, $index, $text) {
    $p = $doc.Paragraphs.Item($index)
    $s = $p.Range.Start
    $e = $p.Range.End
    $r = $doc.Range($s, $e)
    $r.Text = $text
}

# ------------------------------------------------------------------
# 1) Insert "check your city title" (ilvl 1) right after "fix table css"
# ------------------------------------------------------------------
$idxFixTableCss = Find-ParaIndex $d "fix table css"
$pFixTableCss = $d.Paragraphs.Item($idxFixTableCss)
$pFixTableCss.Range.InsertParagraphAfter()
$idxCheckCity = $idxFixTableCss + 1
Set-ParaText $d $idxCheckCity "check your city title"
# new paragraph already inherits ilvl 1 from "fix table css" - leave as-is

# ------------------------------------------------------------------
# 2) Replace the 'use "...interpolateGnBu..."' paragraph (3 runs) with a
#    single run "line - labels"
# ------------------------------------------------------------------
$idxUseQuote = Find-ParaIndex $d 'use "var colorIntrepulate = d3.interpolateGnBu;"'
Set-ParaText $d $idxUseQuote "line - labels"

# ------------------------------------------------------------------
# 3) Replace the big "Main app js : ..." paragraph with "Add short user manual"
# ------------------------------------------------------------------
$idxMainAppJs = Find-ParaIndex $d 'Main app js : svg.select(''#y-axis-text'').text(showPercent ? "מצביעים באחוזים" : "מצביעים");'
Set-ParaText $d $idxMainAppJs "Add short user manual"

# ------------------------------------------------------------------
# 4) Replace the (now duplicate) "Add short user manual" paragraph (the one
#    that originally held this text, further down) with "Total number of
#    votes on each view"
# ------------------------------------------------------------------
$idxOldAddManual = $idxMainAppJs + 1
Set-ParaText $d $idxOldAddManual "Total number of votes on each view"

# ------------------------------------------------------------------
# 5) Insert "Also, sum each cluster's population" (ilvl 1) right after that
# ------------------------------------------------------------------
$pTotalVotes = $d.Paragraphs.Item($idxOldAddManual)
$pTotalVotes.Range.InsertParagraphAfter()
$idxAlsoSum = $idxOldAddManual + 1
Set-ParaText $d $idxAlsoSum "Also, sum each cluster's population"
$pAlsoSum = $d.Paragraphs.Item($idxAlsoSum)
$pAlsoSum.Range.ListFormat.ListLevelNumber = 2

# ------------------------------------------------------------------
# 6) The final (empty) paragraph becomes "Clickable ENGLISH-HEBREW" + "-GERMAN"
#    (two runs) and keeps/receives the _GoBack bookmark.
# ------------------------------------------------------------------
$idxLast = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($idxLast)
Set-ParaText $d $idxLast "Clickable ENGLISH-HEBREW"
$pLast = $d.Paragraphs.Item($idxLast)
$endOfLast = $pLast.Range.End - 1
$rEnd = $d.Range($endOfLast, $endOfLast)
$rEnd.InsertAfter("-GERMAN")

$pLast = $d.Paragraphs.Item($idxLast)
$d.Bookmarks.Add("_GoBack", $pLast.Range)

Write-Output "done"
